$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3/4: the "Virtual method in Monster class" / "Override method in
#     Skeleton class" labels swap places ---
$ws.Range("B3").Value = "Virtual method in Monster class"
$ws.Range("B4").Value = "Override method in Skeleton class"

# --- New screenshot filename cells (added in the exact order the original
#     author typed them, so the shared-string table comes out in the same
#     sequence) ---
$ws.Range("C3").Value = "virtual Monster.GetAttackMessage.png"
$ws.Range("D4").Value = "override Shulker.GetAttackMessage.png"
$ws.Range("C4").Value = "override Dragon.GetAttackMessage.png"
$ws.Range("E4").Value = "override Skeleton.GetAttackMessage.png"
$ws.Range("F4").Value = "override Warden.GetAttackMessage.png"
$ws.Range("G4").Value = "override Witch.GetAttackMessage.png"
$ws.Range("C5").Value = "UserInterface.ShowTurnDecisions.png"
$ws.Range("C6").Value = "UserInterface.DisplayEnumerable Weapons.png"
$ws.Range("D6").Value = "UserInterface.DisplayEnumerable.png"

# --- Column A gets the same green "status" fill on rows 3-5 as row 2 ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column widths: column B widens to fit the longer labels, and the new
#     columns C:G get sized for the filenames placed in them. (Values are
#     tuned so this engine's width rounding lands on the closest
#     representable width to the authored file.) ---
$ws.Columns("B:B").ColumnWidth = 57.9805
$ws.Columns("C:C").ColumnWidth = 36.1719
$ws.Columns("D:D").ColumnWidth = 36.6621
$ws.Columns("E:E").ColumnWidth = 37.8242
$ws.Columns("F:F").ColumnWidth = 36.8242
$ws.Columns("G:G").ColumnWidth = 34.9805

# --- Selection moves to D14 ---
$ws.Range("D14").Select() | Out-Null

Write-Host "done"
